$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 / 18 data first (order matters for shared-string table order) ---
$ws.Range("C19").Value = "SpawnConfig"
$ws.Range("C18").Value = "FootHold"
$ws.Range("E18").Value = "FootHold.xlsx"
$ws.Range("E19").Value = "SpawnConfig.xlsx"
$ws.Range("B18").Value = "FootHoldCategory"
$ws.Range("B19").Value = "SpawnConfigCategory"

# --- Row 17 (BattleLevelConfig...) ---
$b17 = $ws.Range("B17")
$b17.Value = "BattleLevelConfigCategory"
$b17.Font.Family = 3
$b17.Characters(18, 8).Font.Bold = $false

$c17 = $ws.Range("C17")
$c17.Value = "BattleLevelConfig"

$e17 = $ws.Range("E17")
$e17.Value = "BattleLevelConfig.xlsx"
$e17.Font.Family = 3
$e17.Characters(18, 5).Font.Bold = $false

# --- Booleans ---
$ws.Range("D17").Value = $true
$ws.Range("D18").Value = $true
$ws.Range("D19").Value = $true

# --- Re-apply the "category" cell formatting (style index matching B16/E16) ---
# so the cell-level style (s=) matches the rest of the table, without disturbing
# the shared-string text/rich-run content already written above.
$ws.Range("B16").Copy()
$b17.PasteSpecial(-4122)

$ws.Range("B16").Copy()
$e17.PasteSpecial(-4122)
# E17 must end up on its own distinct style (new font family), restore that font.
$e17.Font.Family = 3

$excel.CutCopyMode = 0

# --- Finish up: selection + active cell, like the source edit ---
$ws.Range("B18").Select()
